# Auto-generated PowerShell Excel COM-interop script
# Applies updated numeric values to Sheet1 (pl_mw.xlsx case data)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.3020198183113507
$ws.Cells.Item(2, 3).Value = 0.07458860217634822
$ws.Cells.Item(2, 5).Value = 0.4210905462599328
$ws.Cells.Item(2, 6).Value = 0.4443680307746121
$ws.Cells.Item(2, 7).Value = 0.2923424507497643
$ws.Cells.Item(2, 8).Value = 0.4720819551906956
$ws.Cells.Item(2, 11).Value = 0.3005669203260766
$ws.Cells.Item(2, 14).Value = 1.075236111150243
$ws.Cells.Item(2, 15).Value = 1.454887358339022
$ws.Cells.Item(3, 2).Value = 0.2648994802753464
$ws.Cells.Item(3, 3).Value = 0.07064570305399798
$ws.Cells.Item(3, 5).Value = 0.3674713771598022
$ws.Cells.Item(3, 6).Value = 0.3878228170618172
$ws.Cells.Item(3, 7).Value = 0.2941176245743691
$ws.Cells.Item(3, 8).Value = 0.4763030475964385
$ws.Cells.Item(3, 11).Value = 0.2622639032088614
$ws.Cells.Item(3, 14).Value = 1.085737194912184
$ws.Cells.Item(3, 15).Value = 1.467262182153661
$ws.Cells.Item(4, 2).Value = 0.2420733955499941
$ws.Cells.Item(4, 3).Value = 0.06820961911171253
$ws.Cells.Item(4, 5).Value = 0.334621184651084
$ws.Cells.Item(4, 6).Value = 0.3531389305169483
$ws.Cells.Item(4, 7).Value = 0.2954647070561549
$ws.Cells.Item(4, 8).Value = 0.4791227071872299
$ws.Cells.Item(4, 11).Value = 0.2386568866140664
$ws.Cells.Item(4, 14).Value = 1.092593874484752
$ws.Cells.Item(4, 15).Value = 1.475874947344238
$ws.Cells.Item(5, 2).Value = 0.2327635102925001
$ws.Cells.Item(5, 3).Value = 0.06721314575372617
$ws.Cells.Item(5, 5).Value = 0.3212512965448724
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.2960781256388927
$ws.Cells.Item(5, 8).Value = 0.4803290330666599
$ws.Cells.Item(5, 11).Value = 0.2290149875692862
$ws.Cells.Item(5, 14).Value = 1.095490912002852
$ws.Cells.Item(5, 15).Value = 1.479639420149155
$ws.Cells.Item(6, 2).Value = 0.2312171379869028
$ws.Cells.Item(6, 3).Value = 0.06704745727346051
$ws.Cells.Item(6, 5).Value = 0.3190322149028759
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.2961838718936534
$ws.Cells.Item(6, 8).Value = 0.4805328024517053
$ws.Cells.Item(6, 11).Value = 0.2274126522766693
$ws.Cells.Item(6, 14).Value = 1.095978178266808
$ws.Cells.Item(6, 15).Value = 1.480279878990203
$ws.Cells.Item(7, 2).Value = 0.2419478712071452
$ws.Cells.Item(7, 3).Value = 0.06819619543107081
$ws.Cells.Item(7, 5).Value = 0.3344408072489102
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.2954727190346063
$ws.Cells.Item(7, 8).Value = 0.4791387441672654
$ws.Cells.Item(7, 11).Value = 0.2385269404116173
$ws.Cells.Item(7, 14).Value = 1.092632528308428
$ws.Cells.Item(7, 15).Value = 1.475924685628939
$ws.Cells.Item(8, 2).Value = 0.2892280753377179
$ws.Cells.Item(8, 3).Value = 0.07323226195701693
$ws.Cells.Item(8, 5).Value = 0.4025868613115904
$ws.Cells.Item(8, 6).Value = 0.4248636149813478
$ws.Cells.Item(8, 7).Value = 0.2929010550697839
$ws.Cells.Item(8, 8).Value = 0.4734900807079399
$ws.Cells.Item(8, 11).Value = 0.2873787509428212
$ws.Cells.Item(8, 14).Value = 1.078772030225267
$ws.Cells.Item(8, 15).Value = 1.458943367936385
$ws.Cells.Item(9, 2).Value = 0.38165714854415
$ws.Cells.Item(9, 3).Value = 0.08298599006671736
$ws.Cells.Item(9, 5).Value = 0.5368597020203794
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.2899066807842132
$ws.Cells.Item(9, 8).Value = 0.4642219526549525
$ws.Cells.Item(9, 11).Value = 0.382455577142224
$ws.Cells.Item(9, 14).Value = 1.054834485054933
$ws.Cells.Item(9, 15).Value = 1.43371295806125
$ws.Cells.Item(10, 2).Value = 0.4493732371662134
$ws.Cells.Item(10, 3).Value = 0.09007573695409121
$ws.Cells.Item(10, 5).Value = 0.6359993083288771
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.2889678287747941
$ws.Cells.Item(10, 8).Value = 0.4585164034533804
$ws.Cells.Item(10, 11).Value = 0.4518528070528873
$ws.Cells.Item(10, 14).Value = 1.039221373861032
$ws.Cells.Item(10, 15).Value = 1.420124107684956
$ws.Cells.Item(11, 2).Value = 0.4801343672444318
$ws.Cells.Item(11, 3).Value = 0.09328410422605771
$ws.Cells.Item(11, 5).Value = 0.6812288779715914
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.2888172143815098
$ws.Cells.Item(11, 8).Value = 0.4561607065184035
$ws.Cells.Item(11, 11).Value = 0.483321486305158
$ws.Cells.Item(11, 14).Value = 1.032546547403513
$ws.Cells.Item(11, 15).Value = 1.415022781413981
$ws.Cells.Item(12, 2).Value = 0.4917761771157245
$ws.Cells.Item(12, 3).Value = 0.09449656842126331
$ws.Cells.Item(12, 5).Value = 0.6983765020532786
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.2888001466343297
$ws.Cells.Item(12, 8).Value = 0.4553031694837983
$ws.Cells.Item(12, 11).Value = 0.4952230083470681
$ws.Cells.Item(12, 14).Value = 1.030080432311443
$ws.Cells.Item(12, 15).Value = 1.413246886941948
$ws.Cells.Item(13, 2).Value = 0.4892692144438229
$ws.Cells.Item(13, 3).Value = 0.09423555358375779
$ws.Cells.Item(13, 5).Value = 0.69468253493622
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.2888020416061252
$ws.Cells.Item(13, 8).Value = 0.4554863198591832
$ws.Cells.Item(13, 11).Value = 0.4926604775686485
$ws.Cells.Item(13, 14).Value = 1.030608819193787
$ws.Cells.Item(13, 15).Value = 1.413622417069348
$ws.Cells.Item(14, 2).Value = 0.4810922845805976
$ws.Cells.Item(14, 3).Value = 0.09338390432721155
$ws.Cells.Item(14, 5).Value = 0.6826392110428685
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.2888150080129037
$ws.Cells.Item(14, 8).Value = 0.4560894645203817
$ws.Cells.Item(14, 11).Value = 0.4843009350533123
$ws.Cells.Item(14, 14).Value = 1.032342426129922
$ws.Cells.Item(14, 15).Value = 1.414873550716933
$ws.Cells.Item(15, 2).Value = 0.4760827821905309
$ws.Cells.Item(15, 3).Value = 0.0928619204949257
$ws.Cells.Item(15, 5).Value = 0.6752649957528263
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.2888281613844299
$ws.Cells.Item(15, 8).Value = 0.456463403815782
$ws.Cells.Item(15, 11).Value = 0.4791785066468606
$ws.Cells.Item(15, 14).Value = 1.033412319200607
$ws.Cells.Item(15, 15).Value = 1.415660220144275
$ws.Cells.Item(16, 2).Value = 0.4473620081505487
$ws.Cells.Item(16, 3).Value = 0.08986571954069689
$ws.Cells.Item(16, 5).Value = 0.6330461949412012
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.2889832530961129
$ws.Cells.Item(16, 8).Value = 0.4586751817326018
$ws.Cells.Item(16, 11).Value = 0.4497941872724311
$ws.Cells.Item(16, 14).Value = 1.039666194965633
$ws.Cells.Item(16, 15).Value = 1.4204792743582
$ws.Cells.Item(17, 2).Value = 0.4297312781010305
$ws.Cells.Item(17, 3).Value = 0.08802330321405805
$ws.Cells.Item(17, 5).Value = 0.6071807416197856
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.2891493578941819
$ws.Cells.Item(17, 8).Value = 0.4600934743442266
$ws.Cells.Item(17, 11).Value = 0.4317417444610214
$ws.Cells.Item(17, 14).Value = 1.043612277485799
$ws.Cells.Item(17, 15).Value = 1.423712670789712
$ws.Cells.Item(18, 2).Value = 0.4195865091536461
$ws.Cells.Item(18, 3).Value = 0.08696201689068062
$ws.Cells.Item(18, 5).Value = 0.5923158123999315
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.2892709074676603
$ws.Cells.Item(18, 8).Value = 0.4609318065037868
$ws.Cells.Item(18, 11).Value = 0.4213490322251516
$ws.Cells.Item(18, 14).Value = 1.045922212571988
$ws.Cells.Item(18, 15).Value = 1.425674087607234
$ws.Cells.Item(19, 2).Value = 0.4161509874703881
$ws.Cells.Item(19, 3).Value = 0.08660241452423634
$ws.Cells.Item(19, 5).Value = 0.5872848546218989
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.2893165229945325
$ws.Cells.Item(19, 8).Value = 0.4612195263992405
$ws.Cells.Item(19, 11).Value = 0.4178286361200776
$ws.Cells.Item(19, 14).Value = 1.04671123039946
$ws.Cells.Item(19, 15).Value = 1.426355632779902
$ws.Cells.Item(20, 2).Value = 0.4316085213986014
$ws.Cells.Item(20, 3).Value = 0.08821959524419754
$ws.Cells.Item(20, 5).Value = 0.6099328909520949
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.289128982244975
$ws.Cells.Item(20, 8).Value = 0.459940158754506
$ws.Cells.Item(20, 11).Value = 0.4336644371180967
$ws.Cells.Item(20, 14).Value = 1.043188044256375
$ws.Cells.Item(20, 15).Value = 1.423357945779102
$ws.Cells.Item(21, 2).Value = 0.4834942349392577
$ws.Cells.Item(21, 3).Value = 0.09363412200123378
$ws.Cells.Item(21, 5).Value = 0.6861760686028617
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.2888101131070613
$ws.Cells.Item(21, 8).Value = 0.4559113693373789
$ws.Cells.Item(21, 11).Value = 0.4867567462544287
$ws.Cells.Item(21, 14).Value = 1.031831554725038
$ws.Cells.Item(21, 15).Value = 1.414501828099588
$ws.Cells.Item(22, 2).Value = 0.5173648803925062
$ws.Cells.Item(22, 3).Value = 0.09715837528710836
$ws.Cells.Item(22, 5).Value = 0.7361236047758553
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.2888347674678045
$ws.Cells.Item(22, 8).Value = 0.4534795139850729
$ws.Cells.Item(22, 11).Value = 0.5213680541465919
$ws.Cells.Item(22, 14).Value = 1.024767885221159
$ws.Cells.Item(22, 15).Value = 1.409622597972728
$ws.Cells.Item(23, 2).Value = 0.4992912878941524
$ws.Cells.Item(23, 3).Value = 0.09527875713354206
$ws.Cells.Item(23, 5).Value = 0.7094543810183609
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.2888002151234303
$ws.Cells.Item(23, 8).Value = 0.454759020773551
$ws.Cells.Item(23, 11).Value = 0.5029035387385932
$ws.Cells.Item(23, 14).Value = 1.028505099190646
$ws.Cells.Item(23, 15).Value = 1.412143411722326
$ws.Cells.Item(24, 2).Value = 0.4307598465364322
$ws.Cells.Item(24, 3).Value = 0.08813085802074738
$ws.Cells.Item(24, 5).Value = 0.6086886273486272
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.2891381129345731
$ws.Cells.Item(24, 8).Value = 0.4600094012854683
$ws.Cells.Item(24, 11).Value = 0.4327952317183019
$ws.Cells.Item(24, 14).Value = 1.043379711511768
$ws.Cells.Item(24, 15).Value = 1.423517997725213
$ws.Cells.Item(25, 2).Value = 0.3566850909993491
$ws.Cells.Item(25, 3).Value = 0.08036061027620178
$ws.Cells.Item(25, 5).Value = 0.5004566615563419
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.2904961555504002
$ws.Cells.Item(25, 8).Value = 0.4665354722069779
$ws.Cells.Item(25, 11).Value = 0.3568136602469849
$ws.Cells.Item(25, 14).Value = 1.060963391261154
$ws.Cells.Item(25, 15).Value = 1.439671594184304
